$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5: "AVINA PRODUCE TIN: xxxxx8949" with same values in C/E and D/F
$rows2to5 = 2,3,4,5
foreach ($r in $rows2to5) {
    $ws.Range("B$r").Value = "AVINA PRODUCE TIN: xxxxx8949"
    $ws.Range("C$r").Value = "546,52"
    $ws.Range("D$r").Value = "72,92"
    $ws.Range("E$r").Value = "546,52"
    $ws.Range("F$r").Value = "72,92"
}

# Row 6: "AVINA PRODUCE" (no TIN suffix), E6 differs (1426,00)
$ws.Range("B6").Value = "AVINA PRODUCE"
$ws.Range("C6").Value = "546,52"
$ws.Range("D6").Value = "72,92"
$ws.Range("E6").Value = "1426,00"
$ws.Range("F6").Value = "72,92"
